$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated staff schedule assignments for B2:AC10 (9 staff x 28 days),
# per [ADR-2614] updated alns requirements.
$rows = @(
    @("DO","M3","M1","M1","M1","A1","M1","DO","M1","M3","M1","M1","A1","M1","DO","M1","M1","M1","M1","A1","M3","M1","M1","M1","M3","M1","A1","DO"),
    @("DO","M1","M2","M1","M3","A1","M2","DO","M1","M1","M3","M2","A2","M1","DO","M3","M2","M1","M1","A2","M1","M1","DO","M3","M1","M1","A1","M2"),
    @("DO","M2","M2","M1","M3","M3","M3","M1","DO","M3","M1","M1","A1","M1","M1","DO","M1","M1","M3","A1","M3","M3","M1","M1","DO","M1","A1","M1"),
    @("M1","DO","M1","M3","M2","M2","A2","M2","DO","M1","M3","M1","M2","A1","M1","M2","M1","M1","M3","M2","DO","M1","DO","M1","M3","M2","M2","A1"),
    @("DO","M3","A1","A1","A1","M1","A2","DO","M3","A2","A2","A2","M2","A1","DO","M3","A2","A2","A1","M2","A2","DO","M3","A1","A1","A2","M2","A2"),
    @("DO","M3","A2","A2","A1","A2","A2","M3","A1","A1","A1","A1","DO","A1","M3","A1","DO","A1","A1","A1","A1","M3","A1","A2","A1","A2","A2","DO"),
    @("M3","A2","A1","A2","A1","M2","DO","M3","A2","A2","A1","A1","M1","DO","M3","A2","A2","A1","A2","M1","DO","M3","A2","A1","A1","A2","M2","DO"),
    @("M1","DO","M2","M1","M3","A2","A2","M2","M2","DO","M3","M1","A1","A1","M1","DO","M1","M3","M2","A2","A1","M1","M1","DO","M1","M3","A2","A2"),
    @("A2","M3","DO","M1","M2","M1","M1","M3","M2","DO","M1","M2","M1","M2","M3","M2","DO","A1","A2","A2","M2","DO","M1","A2","M2","M3","A2","M1")
)

$data = New-Object 'object[,]' $rows.Count, $rows[0].Count
for ($i = 0; $i -lt $rows.Count; $i++) {
    for ($j = 0; $j -lt $rows[$i].Count; $j++) {
        $data[$i, $j] = $rows[$i][$j]
    }
}

$ws.Range("B2:AC10").Value = $data
